$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume refresh. Price values in column D are stored
# as text in the source data (e.g. "239.09"), so for values that Excel would
# otherwise auto-convert to a number, force the Text number format before
# assigning the value so the literal string is preserved.

$ws.Range('D2').Value = '35.200.81'
$ws.Range('E2').Value = '  +1.22%  '
$ws.Range('D3').Value = '1.858.57'
$ws.Range('E3').Value = '  +1.71%  '
$ws.Range('E4').Value = '  +0.54%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '239.09'
$ws.Range('E5').Value = '  +3.84%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.623'
$ws.Range('E6').Value = '  +0.78%  '
$ws.Range('E7').Value = '  +0.47%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '42.27'
$ws.Range('E8').Value = '  +6.83%  '
$ws.Range('E9').Value = '  +1.21%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0693'
$ws.Range('E10').Value = '  +1.67%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0988'
$ws.Range('E11').Value = '  +0.14%  '
$ws.Range('E12').Value = '  +1.58%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.51'
$ws.Range('E13').Value = '  +1.81%  '
$ws.Range('D14').Value = '1.870.96'
$ws.Range('E14').Value = '  +2.04%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.677'
$ws.Range('E15').Value = '  +1.33%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '4.73'
$ws.Range('E16').Value = '  +2.36%  '
$ws.Range('D17').Value = '35.164.05'
$ws.Range('E17').Value = '  +1.10%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '69.89'
$ws.Range('E18').Value = '  +0.62%  '
$ws.Range('E19').Value = '  +1.35%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '240.85'
$ws.Range('E20').Value = '  +0.34%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.22'
$ws.Range('E21').Value = '  +0.85%  '
$ws.Range('E22').Value = '  +1.99%  '
$ws.Range('E23').Value = '  +0.38%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.28'
$ws.Range('E24').Value = '  +1.23%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '169.91'
$ws.Range('E25').Value = '  -1.04%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.90'
$ws.Range('E26').Value = '  +26.60%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.02'
$ws.Range('E27').Value = '  +3.65%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '17.66'
$ws.Range('E28').Value = '  +2.08%  '
$ws.Range('E29').Value = '  +0.15%  '
$ws.Range('B30').Value = 'Hedera'
$ws.Range('C30').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0560'
$ws.Range('E30').Value = '  +1.95%  '
$ws.Range('B31').Value = 'BinanceUSD'
$ws.Range('C31').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.01'
$ws.Range('E31').Value = '  +0.46%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.00'
$ws.Range('E32').Value = '  +2.16%  '
$ws.Range('E33').Value = '  +28.11%  '
$ws.Range('E34').Value = '  +2.32%  '
$ws.Range('E35').Value = '  +10.71%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.813'
$ws.Range('E36').Value = '  +16.79%  '
$ws.Range('E37').Value = '  +8.03%  '
$ws.Range('E38').Value = '  +4.60%  '
$ws.Range('E39').Value = '  +4.56%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '90.06'
$ws.Range('E40').Value = '  -1.05%  '
$ws.Range('D41').Value = '1.347.30'
$ws.Range('E41').Value = '  +0.76%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0600'
$ws.Range('E42').Value = '  +15.11%  '
$ws.Range('E43').Value = '  +3.46%  '
$ws.Range('E44').Value = '  +3.06%  '
$ws.Range('E45').Value = '  +0.38%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '12.45'
$ws.Range('E46').Value = '  +40.47%  '
$ws.Range('E47').Value = '  -0.50%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '6.58'
$ws.Range('E48').Value = '  +5.46%  '
$ws.Range('E49').Value = '  +1.89%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0681'
$ws.Range('E50').Value = '  +1.31%  '
$ws.Range('E51').Value = '  +0.52%  '
